$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.692.02'
$ws.Range("E2").Value = '  -1.49%  '
$ws.Range("D3").Value = '3.266.45'
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("E4").Value = '  -0.03%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '580.08'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.58%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '184.13'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.27%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +0.84%  '
$ws.Range("E9").Value = '  -2.74%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '6.57'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -1.39%  '
$ws.Range("E11").Value = '  -4.55%  '
$ws.Range("D12").Value = '3.831.15'
$ws.Range("E12").Value = '  -0.52%  '
$ws.Range("E13").Value = '  +0.57%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '27.33'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -4.67%  '
$ws.Range("D15").Value = '67.714.44'
$ws.Range("E15").Value = '  -1.44%  '
$ws.Range("E16").Value = '  -2.05%  '
$ws.Range("D17").Value = '3.264.53'
$ws.Range("E17").Value = '  -0.76%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '5.69'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -2.45%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '13.40'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.16%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '401.02'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +1.65%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '7.54'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.95%  '
$ws.Range("E22").Value = '  +0.12%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '71.01'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.48%  '
$ws.Range("E24").Value = '  -1.68%  '
$ws.Range("E25").Value = '  -1.75%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.187'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -1.42%  '
$ws.Range("E27").Value = '  -2.09%  '
$ws.Range("E28").Value = '  +0.44%  '
$ws.Range("E29").Value = '  -1.71%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '22.61'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -1.43%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '5.46'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -4.10%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '6.90'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -3.19%  '
$ws.Range("E34").Value = '  -4.09%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '164.16'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.25%  '
$ws.Range("E36").Value = '  -3.75%  '
$ws.Range("E37").Value = '  -2.34%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '26.99'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +2.71%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.802'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -3.47%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '4.49'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.98%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '6.35'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -3.50%  '
$ws.Range("D42").Value = '2.679.80'
$ws.Range("E42").Value = '  +2.21%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '40.70'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -1.70%  '
$ws.Range("E44").Value = '  -3.94%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0676'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -1.51%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '334.56'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -2.68%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '24.52'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.07%  '
$ws.Range("E48").Value = '  -3.23%  '
$ws.Range("E49").Value = '  -0.65%  '
$ws.Range("E50").Value = '  -1.82%  '
$ws.Range("B51").Value = 'FirstDigitalUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.05%  '
